$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 937: correct the Event text (G column) to the variant with a leading space ---
$ws.Cells.Item(937, 7).Value = ' Girl kisses boy'

# --- Header: add new "Comments" column L ---
$ws.Cells.Item(1, 12).Value = "Comments"

# --- Fill Agent/Verb/Patient Recode columns (I, J, K) for newly blind-coded rows, plus two analyst comments (L) ---
$ws.Cells.Item(890, 9).Value = 1
$ws.Cells.Item(890, 10).Value = 1
$ws.Cells.Item(890, 11).Value = 1
$ws.Cells.Item(891, 9).Value = 1
$ws.Cells.Item(891, 10).Value = 0
$ws.Cells.Item(891, 11).Value = 0
$ws.Cells.Item(892, 9).Value = 1
$ws.Cells.Item(892, 10).Value = 1
$ws.Cells.Item(892, 11).Value = 1
$ws.Cells.Item(893, 9).Value = 1
$ws.Cells.Item(893, 10).Value = 1
$ws.Cells.Item(893, 11).Value = 0
$ws.Cells.Item(894, 9).Value = 1
$ws.Cells.Item(914, 9).Value = 1
$ws.Cells.Item(914, 10).Value = 1
$ws.Cells.Item(914, 11).Value = 0
$ws.Cells.Item(915, 9).Value = 1
$ws.Cells.Item(915, 10).Value = 1
$ws.Cells.Item(915, 11).Value = 0
$ws.Cells.Item(916, 9).Value = 1
$ws.Cells.Item(916, 10).Value = 0
$ws.Cells.Item(916, 11).Value = 0
$ws.Cells.Item(917, 9).Value = 1
$ws.Cells.Item(917, 10).Value = 1
$ws.Cells.Item(917, 11).Value = 1
$ws.Cells.Item(918, 9).Value = 1
$ws.Cells.Item(918, 10).Value = 1
$ws.Cells.Item(918, 11).Value = 0
$ws.Cells.Item(919, 9).Value = 0
$ws.Cells.Item(919, 10).Value = 0
$ws.Cells.Item(919, 11).Value = 0
$ws.Cells.Item(920, 9).Value = 1
$ws.Cells.Item(920, 10).Value = 1
$ws.Cells.Item(920, 11).Value = 1
$ws.Cells.Item(921, 9).Value = 1
$ws.Cells.Item(921, 10).Value = 1
$ws.Cells.Item(921, 11).Value = 1
$ws.Cells.Item(922, 9).Value = 1
$ws.Cells.Item(922, 10).Value = 1
$ws.Cells.Item(922, 11).Value = 0
$ws.Cells.Item(923, 9).Value = 0
$ws.Cells.Item(923, 10).Value = 0
$ws.Cells.Item(923, 11).Value = 0
$ws.Cells.Item(924, 9).Value = 1
$ws.Cells.Item(924, 10).Value = 1
$ws.Cells.Item(924, 11).Value = 1
$ws.Cells.Item(925, 9).Value = 1
$ws.Cells.Item(925, 10).Value = 0
$ws.Cells.Item(925, 11).Value = 0
$ws.Cells.Item(926, 9).Value = 1
$ws.Cells.Item(926, 10).Value = 1
$ws.Cells.Item(926, 11).Value = 1
$ws.Cells.Item(927, 9).Value = 1
$ws.Cells.Item(927, 10).Value = 1
$ws.Cells.Item(927, 11).Value = 0
$ws.Cells.Item(928, 9).Value = 1
$ws.Cells.Item(928, 10).Value = 1
$ws.Cells.Item(928, 11).Value = 1
$ws.Cells.Item(929, 9).Value = 1
$ws.Cells.Item(929, 10).Value = 1
$ws.Cells.Item(929, 11).Value = 1
$ws.Cells.Item(930, 9).Value = 0
$ws.Cells.Item(930, 10).Value = 0
$ws.Cells.Item(930, 11).Value = 0
$ws.Cells.Item(931, 9).Value = 1
$ws.Cells.Item(931, 10).Value = 1
$ws.Cells.Item(931, 11).Value = 0
$ws.Cells.Item(932, 9).Value = 1
$ws.Cells.Item(932, 10).Value = 0
$ws.Cells.Item(932, 11).Value = 0
$ws.Cells.Item(933, 9).Value = 1
$ws.Cells.Item(933, 10).Value = 1
$ws.Cells.Item(933, 11).Value = 1
$ws.Cells.Item(934, 9).Value = 1
$ws.Cells.Item(934, 10).Value = 1
$ws.Cells.Item(934, 11).Value = 0
$ws.Cells.Item(935, 9).Value = 1
$ws.Cells.Item(935, 10).Value = 1
$ws.Cells.Item(935, 11).Value = 1
$ws.Cells.Item(936, 9).Value = 1
$ws.Cells.Item(936, 10).Value = 0
$ws.Cells.Item(936, 11).Value = 0
$ws.Cells.Item(937, 9).Value = 1
$ws.Cells.Item(937, 10).Value = 1
$ws.Cells.Item(937, 11).Value = 0
$ws.Cells.Item(937, 12).Value = 'Look at video 242 and 248 where he gestures "boy lifts girl". He always does boy and girl the same way which leads me to believe the word order on this video 246 "girl kisses boy" is SOV. He hesitates after doing the gesture for boy and then does girl. '
$ws.Cells.Item(938, 9).Value = 1
$ws.Cells.Item(938, 10).Value = 0
$ws.Cells.Item(938, 11).Value = 0
$ws.Cells.Item(939, 9).Value = 1
$ws.Cells.Item(939, 10).Value = 0
$ws.Cells.Item(939, 11).Value = 0
$ws.Cells.Item(940, 9).Value = 1
$ws.Cells.Item(940, 10).Value = 0
$ws.Cells.Item(940, 11).Value = 0
$ws.Cells.Item(941, 9).Value = 1
$ws.Cells.Item(941, 10).Value = 0
$ws.Cells.Item(941, 11).Value = 1
$ws.Cells.Item(942, 9).Value = 1
$ws.Cells.Item(942, 10).Value = 0
$ws.Cells.Item(942, 11).Value = 0
$ws.Cells.Item(943, 9).Value = 0
$ws.Cells.Item(943, 10).Value = 0
$ws.Cells.Item(943, 11).Value = 0
$ws.Cells.Item(944, 9).Value = 1
$ws.Cells.Item(944, 10).Value = 0
$ws.Cells.Item(944, 11).Value = 0
$ws.Cells.Item(945, 9).Value = 1
$ws.Cells.Item(945, 10).Value = 0
$ws.Cells.Item(945, 11).Value = 1
$ws.Cells.Item(946, 9).Value = 1
$ws.Cells.Item(946, 10).Value = 0
$ws.Cells.Item(946, 11).Value = 0
$ws.Cells.Item(947, 9).Value = 0
$ws.Cells.Item(947, 10).Value = 0
$ws.Cells.Item(947, 11).Value = 0
$ws.Cells.Item(948, 9).Value = 1
$ws.Cells.Item(948, 10).Value = 0
$ws.Cells.Item(948, 11).Value = 1
$ws.Cells.Item(949, 9).Value = 1
$ws.Cells.Item(949, 10).Value = 0
$ws.Cells.Item(949, 11).Value = 0
$ws.Cells.Item(950, 9).Value = 1
$ws.Cells.Item(950, 10).Value = 0
$ws.Cells.Item(950, 11).Value = 1
$ws.Cells.Item(951, 9).Value = 1
$ws.Cells.Item(951, 10).Value = 0
$ws.Cells.Item(951, 11).Value = 0
$ws.Cells.Item(952, 9).Value = 1
$ws.Cells.Item(952, 10).Value = 1
$ws.Cells.Item(952, 11).Value = 1
$ws.Cells.Item(953, 9).Value = 1
$ws.Cells.Item(953, 10).Value = 0
$ws.Cells.Item(953, 11).Value = 1
$ws.Cells.Item(954, 9).Value = 0
$ws.Cells.Item(954, 10).Value = 0
$ws.Cells.Item(954, 11).Value = 0
$ws.Cells.Item(955, 9).Value = 1
$ws.Cells.Item(955, 10).Value = 0
$ws.Cells.Item(955, 11).Value = 0
$ws.Cells.Item(956, 9).Value = 0
$ws.Cells.Item(956, 10).Value = 0
$ws.Cells.Item(956, 11).Value = 0
$ws.Cells.Item(957, 9).Value = 1
$ws.Cells.Item(957, 10).Value = 0
$ws.Cells.Item(957, 11).Value = 0
$ws.Cells.Item(958, 9).Value = 1
$ws.Cells.Item(958, 10).Value = 0
$ws.Cells.Item(958, 11).Value = 0
$ws.Cells.Item(959, 9).Value = 1
$ws.Cells.Item(959, 10).Value = 0
$ws.Cells.Item(959, 11).Value = 1
$ws.Cells.Item(960, 9).Value = 1
$ws.Cells.Item(960, 10).Value = 0
$ws.Cells.Item(960, 11).Value = 0
$ws.Cells.Item(961, 9).Value = 1
$ws.Cells.Item(961, 10).Value = 0
$ws.Cells.Item(961, 11).Value = 1
$ws.Cells.Item(962, 9).Value = 0
$ws.Cells.Item(962, 10).Value = 1
$ws.Cells.Item(962, 11).Value = 0
$ws.Cells.Item(963, 9).Value = 0
$ws.Cells.Item(963, 10).Value = 1
$ws.Cells.Item(963, 11).Value = 0
$ws.Cells.Item(964, 9).Value = 1
$ws.Cells.Item(964, 10).Value = 1
$ws.Cells.Item(964, 11).Value = 0
$ws.Cells.Item(965, 9).Value = 1
$ws.Cells.Item(965, 10).Value = 1
$ws.Cells.Item(965, 11).Value = 1
$ws.Cells.Item(966, 9).Value = 1
$ws.Cells.Item(966, 10).Value = 1
$ws.Cells.Item(966, 11).Value = 0
$ws.Cells.Item(967, 9).Value = 0
$ws.Cells.Item(967, 10).Value = 0
$ws.Cells.Item(967, 11).Value = 0
$ws.Cells.Item(968, 9).Value = 1
$ws.Cells.Item(968, 10).Value = 1
$ws.Cells.Item(968, 11).Value = 0
$ws.Cells.Item(969, 9).Value = 1
$ws.Cells.Item(969, 10).Value = 1
$ws.Cells.Item(969, 11).Value = 0
$ws.Cells.Item(970, 9).Value = 1
$ws.Cells.Item(970, 10).Value = 1
$ws.Cells.Item(970, 11).Value = 0
$ws.Cells.Item(971, 9).Value = 0
$ws.Cells.Item(971, 10).Value = 0
$ws.Cells.Item(971, 11).Value = 0
$ws.Cells.Item(972, 9).Value = 0
$ws.Cells.Item(972, 10).Value = 1
$ws.Cells.Item(972, 11).Value = 0
$ws.Cells.Item(972, 12).Value = 'Watch in previous clips how he does old lady (always points to glasses)  and watch how he does girl in clip 256 "Fireman kicks girl. Based on that, in this clip 256 "girl elbows old lady" the word order might be OSV. Old lady  when he points to glasses and girl might be the triangle figure he draws. '
$ws.Cells.Item(973, 9).Value = 1
$ws.Cells.Item(973, 10).Value = 0
$ws.Cells.Item(973, 11).Value = 0
$ws.Cells.Item(974, 9).Value = 1
$ws.Cells.Item(974, 10).Value = 1
$ws.Cells.Item(974, 11).Value = 0
$ws.Cells.Item(975, 9).Value = 1
$ws.Cells.Item(975, 10).Value = 1
$ws.Cells.Item(975, 11).Value = 0
$ws.Cells.Item(976, 9).Value = 1
$ws.Cells.Item(976, 10).Value = 1
$ws.Cells.Item(976, 11).Value = 1
$ws.Cells.Item(977, 9).Value = 0
$ws.Cells.Item(977, 10).Value = 1
$ws.Cells.Item(977, 11).Value = 1
$ws.Cells.Item(978, 9).Value = 0
$ws.Cells.Item(978, 10).Value = 0
$ws.Cells.Item(978, 11).Value = 0
$ws.Cells.Item(979, 9).Value = 1
$ws.Cells.Item(979, 10).Value = 0
$ws.Cells.Item(979, 11).Value = 0
$ws.Cells.Item(980, 9).Value = 0
$ws.Cells.Item(980, 10).Value = 0
$ws.Cells.Item(980, 11).Value = 0
$ws.Cells.Item(981, 9).Value = 1
$ws.Cells.Item(981, 10).Value = 1
$ws.Cells.Item(981, 11).Value = 0
$ws.Cells.Item(982, 9).Value = 1
$ws.Cells.Item(982, 10).Value = 1
$ws.Cells.Item(982, 11).Value = 0
$ws.Cells.Item(983, 9).Value = 1
$ws.Cells.Item(983, 10).Value = 1
$ws.Cells.Item(983, 11).Value = 1
$ws.Cells.Item(984, 9).Value = 0
$ws.Cells.Item(984, 10).Value = 0
$ws.Cells.Item(984, 11).Value = 0
$ws.Cells.Item(985, 9).Value = 0
$ws.Cells.Item(985, 10).Value = 1
$ws.Cells.Item(985, 11).Value = 1
$ws.Cells.Item(986, 9).Value = 1
$ws.Cells.Item(986, 10).Value = 0
$ws.Cells.Item(986, 11).Value = 0
$ws.Cells.Item(987, 9).Value = 0
$ws.Cells.Item(987, 10).Value = 0
$ws.Cells.Item(987, 11).Value = 0
$ws.Cells.Item(988, 9).Value = 1
$ws.Cells.Item(988, 10).Value = 0
$ws.Cells.Item(988, 11).Value = 0
$ws.Cells.Item(989, 9).Value = 1
$ws.Cells.Item(989, 10).Value = 0
$ws.Cells.Item(989, 11).Value = 1
$ws.Cells.Item(990, 9).Value = 1
$ws.Cells.Item(990, 10).Value = 0
$ws.Cells.Item(990, 11).Value = 0
$ws.Cells.Item(991, 9).Value = 0
$ws.Cells.Item(991, 10).Value = 0
$ws.Cells.Item(991, 11).Value = 0
$ws.Cells.Item(992, 9).Value = 0
$ws.Cells.Item(992, 10).Value = 0
$ws.Cells.Item(992, 11).Value = 0
$ws.Cells.Item(993, 9).Value = 1
$ws.Cells.Item(993, 10).Value = 0
$ws.Cells.Item(993, 11).Value = 0
$ws.Cells.Item(994, 9).Value = 1
$ws.Cells.Item(994, 10).Value = 0
$ws.Cells.Item(994, 11).Value = 0
$ws.Cells.Item(995, 9).Value = 0
$ws.Cells.Item(995, 10).Value = 0
$ws.Cells.Item(995, 11).Value = 0
$ws.Cells.Item(996, 9).Value = 0
$ws.Cells.Item(996, 10).Value = 0
$ws.Cells.Item(996, 11).Value = 1
$ws.Cells.Item(997, 9).Value = 1
$ws.Cells.Item(997, 10).Value = 0
$ws.Cells.Item(997, 11).Value = 0

# --- Restore the active selection to match where the coder ended up ---
$ws.Range("I1004").Select()
